# Refresh the cryptos price list (GitHub Actions scrape update).
# Price/volume figures are refreshed for most rows; a few coins also
# swapped rank position (and therefore row) with their neighbour.
# NumberFormat is forced to "@" (Text) right before writing any Price
# value that would otherwise be auto-parsed as a number by Excel, so
# that values like "580.50" or "1.00" are kept as literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.884.39"
$ws.Range("E2").Value = "  +2.13%  "
$ws.Range("D3").Value = "3.100.70"
$ws.Range("E3").Value = "  +5.31%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.50"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.96"
$ws.Range("E6").Value = "  +6.48%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.095.53"
$ws.Range("E8").Value = "  +5.24%  "
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  +3.72%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.44"
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("E12").Value = "  +4.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.38"
$ws.Range("E14").Value = "  +7.53%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "3.612.67"
$ws.Range("E16").Value = "  +5.18%  "
$ws.Range("D17").Value = "66.876.42"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").Value = "3.099.32"
$ws.Range("E19").Value = "  +5.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.21"
$ws.Range("E20").Value = "  +3.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "480.77"
$ws.Range("E21").Value = "  +8.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.715"
$ws.Range("E22").Value = "  +2.87%  "
$ws.Range("E23").Value = "  +3.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.98"
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.17"
$ws.Range("E25").Value = "  +7.85%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("E26").Value = "  +4.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.08"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("E30").Value = "  -3.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.90"
$ws.Range("E32").Value = "  +6.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000101"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  +3.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.991"
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.09"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("E39").Value = "  +7.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.19"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("E41").Value = "  +3.91%  "
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0362"
$ws.Range("E45").Value = "  +3.08%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.840.80"
$ws.Range("E46").Value = "  +6.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "385.85"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.18"
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.02"
$ws.Range("E50").Value = "  +4.29%  "
$ws.Range("E51").Value = "  +2.99%  "
